$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows (A: Oyuncu Adı, B: Pozisyon, C: Takım), rows 2-16 reshuffled;
# rows 17-18 remain unchanged.
$data = @(
    @("Damian Lillard",        "PG",       "Milwaukee Bucks"),
    @("LaMelo Ball",           "PG,SG",    "Charlotte Hornets"),
    @("Brandon Miller",        "SG,SF",    "Charlotte Hornets"),
    @("Anthony Davis",         "PF,C",     "Los Angeles Lakers"),
    @("Robert Williams III",   "C",        "Portland Trail Blazers"),
    @("Bam Adebayo",           "C",        "Miami Heat"),
    @("T.J. McConnell",        "PG",       "Indiana Pacers"),
    @("Isaiah Hartenstein",    "C",        "Oklahoma City Thunder"),
    @("Julius Randle",         "PF",       "Minnesota Timberwolves"),
    @("Bilal Coulibaly",       "SG,SF",    "Washington Wizards"),
    @("Cam Thomas",            "SG,SF",    "Brooklyn Nets"),
    @("Jared McCain",          "PG,SG",    "Philadelphia 76ers"),
    @("Derrick White",         "PG,SG",    "Boston Celtics"),
    @("Brandon Ingram",        "SG,SF,PF", "New Orleans Pelicans"),
    @("Cameron Johnson",       "SF,PF",    "Brooklyn Nets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
